$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: Natrium-Penicillin warning - append sentence about not mixing Penicillin/Aminoglykosides
$ws.Range("I4").Value = "Bad penetration into the central nervous system and abscesses. Commonly administered in combination with Gentamicin. Do not mix Penicillin and Aminoglykosides (Gentamicin/Amikacin) in one syringe."

# Row 6: Gentamicin warning - append sentence
$ws.Range("I6").Value = "Does not work in abscesses or in anaerobic milieu. In Foals <1 month, use Amikacin instead. Do not mix Penicillin and Aminoglykosides (Gentamicin/Amikacin) in one syringe."

# Row 7: Amikacin warning - fix typo "abscesser" -> "abscesses" and append sentence
$ws.Range("I7").Value = "First line antibiotic in foals. Drug repurposing from human medicine. Does not work in abscesses or in anaerobic milieu. Do not mix Penicillin and Aminoglykosides (Gentamicin/Amikacin) in one syringe."

# Row 8 & 9: Oxytetracyclin / Doxycyclin warning - "human medicine" -> "veterinary medicine"
$ws.Range("I8").Value = "Drug repurposing from veterinary medicine, intracellular effect. Works against Ehrlichia, Rickettsia, Anaplasma. Does not penetrate the blood-brain-barrier; in general good tissue penetration. Careful use of Tetracyclines in growing patients."
$ws.Range("I9").Value = "Drug repurposing from veterinary medicine, intracellular effect. Works against Ehrlichia, Rickettsia, Anaplasma. Does not penetrate the blood-brain-barrier; in general good tissue penetration. Careful use of Tetracyclines in growing patients."

# Row 18 & 19: Ceftiofur warning - "human medicine" -> "veterinary medicine"
$ws.Range("I18").Value = "Drug repurposing from veterinary medicine."
$ws.Range("I19").Value = "Drug repurposing from veterinary medicine."

# Update sheet view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G8").Select()
